$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2151364.2
$ws.Range("I38").Value = 4032391
$ws.Range("J38").Value = 1619
$ws.Range("K38").Value = 12097173
$ws.Range("L38").Value = 4857
$ws.Range("M38").Value = -12096801
$ws.Range("N38").Value = -5601
$ws.Range("H58").Value = 330791.47
$ws.Range("I58").Value = 688990.4399999999
$ws.Range("K58").Value = 2066971.32
$ws.Range("M58").Value = -2066821.32
$ws.Range("H62").Value = 1976.5
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 1953
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 1953
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -3201
$ws.Range("H65").Value = 1976.5
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 1953
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 9765
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -16005
$ws.Range("H113").Value = 92655.37
$ws.Range("I113").Value = 126651.125
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 126651.125
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -123397.125
$ws.Range("N113").Value = -8508
$ws.Range("H129").Value = 878.6389
$ws.Range("J129").Value = 889.45715
$ws.Range("L129").Value = 2668.37145
$ws.Range("N129").Value = -12668.37145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 39788.58
$ws.Range("I2").Value = 1158.238
$ws.Range("J2").Value = 202036
$ws.Range("K2").Value = 1158.238
$ws.Range("L2").Value = 202036
$ws.Range("M2").Value = -1045.238
$ws.Range("N2").Value = -202262
$ws.Range("H45").Value = 1747.5
$ws.Range("I45").Value = 1216.2858
$ws.Range("J45").Value = 2278.7144
$ws.Range("K45").Value = 1216.2858
$ws.Range("L45").Value = 2278.7144
$ws.Range("M45").Value = -839.2858000000001
$ws.Range("N45").Value = -3032.7144
$ws.Range("H116").Value = 39788.58
$ws.Range("I116").Value = 1158.238
$ws.Range("J116").Value = 202036
$ws.Range("K116").Value = 1158.238
$ws.Range("L116").Value = 202036
$ws.Range("M116").Value = 1135.762
$ws.Range("N116").Value = -206624
$ws.Range("H123").Value = 42000
$ws.Range("J123").Value = 42000
$ws.Range("L123").Value = 42000
$ws.Range("N123").Value = -51800
$ws.Range("H132").Value = 2463.2642
$ws.Range("I132").Value = 2341.851
$ws.Range("K132").Value = 7025.553
$ws.Range("M132").Value = -4495.553

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 39788.58
$ws.Range("I3").Value = 1158.238
$ws.Range("J3").Value = 202036
$ws.Range("K3").Value = 1158.238
$ws.Range("L3").Value = 202036
$ws.Range("M3").Value = -1044.238
$ws.Range("N3").Value = -202264
$ws.Range("H6").Value = 33000
$ws.Range("J6").Value = 33000
$ws.Range("L6").Value = 33000
$ws.Range("N6").Value = -33226
$ws.Range("H13").Value = 45000
$ws.Range("J13").Value = 45000
$ws.Range("L13").Value = 45000
$ws.Range("N13").Value = -45336
$ws.Range("H62").Value = 49316.668
$ws.Range("I62").Value = 48000
$ws.Range("J62").Value = 49975
$ws.Range("K62").Value = 48000
$ws.Range("L62").Value = 49975
$ws.Range("N62").Value = -51347
$ws.Range("M62").Value = -47314
$ws.Range("H64").Value = 388.14285
$ws.Range("I64").Value = 81
$ws.Range("K64").Value = 81
$ws.Range("M64").Value = 144
$ws.Range("H65").Value = 49316.668
$ws.Range("I65").Value = 48000
$ws.Range("J65").Value = 49975
$ws.Range("K65").Value = 144000
$ws.Range("L65").Value = 149925
$ws.Range("N65").Value = -156789
$ws.Range("M65").Value = -140568
$ws.Range("H67").Value = 388.14285
$ws.Range("I67").Value = 81
$ws.Range("K67").Value = 81
$ws.Range("M67").Value = 699
$ws.Range("H107").Value = 43479204
$ws.Range("I107").Value = 76924100
$ws.Range("J107").Value = 843.5
$ws.Range("K107").Value = 76924100
$ws.Range("L107").Value = 843.5
$ws.Range("M107").Value = -76922180
$ws.Range("N107").Value = -4683.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 19334.5
$ws.Range("J99").Value = 22205.4
$ws.Range("L99").Value = 22205.4
$ws.Range("N99").Value = -25201.4
$ws.Range("H116").Value = 49985
$ws.Range("J116").Value = 49985
$ws.Range("L116").Value = 49985
$ws.Range("N116").Value = -59163
$ws.Range("H126").Value = 19334.5
$ws.Range("J126").Value = 22205.4
$ws.Range("L126").Value = 66616.20000000001
$ws.Range("N126").Value = -71556.20000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 7300
$ws.Range("J87").Value = 12000
$ws.Range("L87").Value = 36000
$ws.Range("N87").Value = -38496
$ws.Range("H90").Value = 7300
$ws.Range("J90").Value = 12000
$ws.Range("L90").Value = 108000
$ws.Range("N90").Value = -120480
$ws.Range("H99").Value = 2304.8
$ws.Range("I99").Value = 1841.3334
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 5524.0002
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = -3278.0002
$ws.Range("N99").Value = -13492
$ws.Range("H113").Value = 564.24
$ws.Range("I113").Value = 606.125
$ws.Range("J113").Value = 544.5294
$ws.Range("K113").Value = 1818.375
$ws.Range("L113").Value = 1633.5882
$ws.Range("M113").Value = 351.625
$ws.Range("N113").Value = -5973.5882
$ws.Range("H131").Value = 778.74
$ws.Range("I131").Value = 408.0909
$ws.Range("J131").Value = 824.55054
$ws.Range("K131").Value = 1224.2727
$ws.Range("L131").Value = 2473.65162
$ws.Range("M131").Value = 3815.7273
$ws.Range("N131").Value = -12553.65162
$ws.Range("H132").Value = 2975.7778
$ws.Range("J132").Value = 4749.75
$ws.Range("L132").Value = 42747.75
$ws.Range("N132").Value = -47807.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 13034.429
$ws.Range("J48").Value = 13034.429
$ws.Range("L48").Value = 13034.429
$ws.Range("N48").Value = -14004.429
$ws.Range("H113").Value = 2206.875
$ws.Range("I113").Value = 3280.2
$ws.Range("J113").Value = 1719
$ws.Range("K113").Value = 3280.2
$ws.Range("L113").Value = 1719
$ws.Range("M113").Value = -1110.2
$ws.Range("N113").Value = -6059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 15854.286
$ws.Range("J43").Value = 15854.286
$ws.Range("L43").Value = 15854.286
$ws.Range("N43").Value = -16240.286
$ws.Range("H117").Value = 41743
$ws.Range("J117").Value = 41743
$ws.Range("L117").Value = 41743
$ws.Range("N117").Value = -50921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8883.857
$ws.Range("I14").Value = 2064.5
$ws.Range("J14").Value = 49800
$ws.Range("K14").Value = 2064.5
$ws.Range("L14").Value = 49800
$ws.Range("M14").Value = -1896.5
$ws.Range("N14").Value = -50136
$ws.Range("H54").Value = 4472.846
$ws.Range("J54").Value = 4256.4165
$ws.Range("L54").Value = 4256.4165
$ws.Range("N54").Value = -5296.4165
$ws.Range("H81").Value = 182958.36
$ws.Range("J81").Value = 252239.75
$ws.Range("L81").Value = 504479.5
$ws.Range("N81").Value = -506601.5
$ws.Range("H84").Value = 182958.36
$ws.Range("J84").Value = 252239.75
$ws.Range("L84").Value = 2522397.5
$ws.Range("N84").Value = -2533005.5
$ws.Range("H118").Value = 39176.4
$ws.Range("J118").Value = 39176.4
$ws.Range("L118").Value = 39176.4
$ws.Range("N118").Value = -42490.4
$ws.Range("H126").Value = 2048.682
$ws.Range("I126").Value = 2087.3333
$ws.Range("K126").Value = 6261.999899999999
$ws.Range("M126").Value = -3791.999899999999
$ws.Range("H136").Value = 883.36365
$ws.Range("I136").Value = 655.6
$ws.Range("J136").Value = 1371.4286
$ws.Range("K136").Value = 1966.8
$ws.Range("L136").Value = 4114.2858
$ws.Range("M136").Value = 583.1999999999998
$ws.Range("N136").Value = -9214.2858
